$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Ripal Patel"
$ws.Name = "Ripal Patel"

# Insert a new column at A (matchNo), shifting existing teamName..result
# columns from A:L to B:M
$ws.Columns.Item(1).Insert()

# Insert a new row after current row 1 (the header), pushing the existing
# data row (Royal Challengers Bangalore match) down from row 2 to row 3
$ws.Rows.Item(2).Insert()

# Make sure the freshly inserted column/row cells are text-formatted too
$ws.Range("A1:M3").NumberFormat = "@"

# New header for the inserted column
$ws.Cells.Item(1, 1).Value = "matchNo"

# Fill in the new row 2 with the "50th" match data (vs Chennai Super Kings)
$ws.Cells.Item(2, 1).Value = "50th"
$ws.Cells.Item(2, 2).Value = "Delhi Capitals"
$ws.Cells.Item(2, 3).Value = "Ripal Patel"
$ws.Cells.Item(2, 4).Value = "c Chahar b Jadeja"
$ws.Cells.Item(2, 5).Value = "18"
$ws.Cells.Item(2, 6).Value = "20"
$ws.Cells.Item(2, 7).Value = "2"
$ws.Cells.Item(2, 8).Value = "0"
$ws.Cells.Item(2, 9).Value = "90.00"
$ws.Cells.Item(2, 10).Value = "Chennai Super Kings"
$ws.Cells.Item(2, 11).Value = "Dubai (DSC)"
$ws.Cells.Item(2, 12).Value = "October 04"
$ws.Cells.Item(2, 13).Value = "Capitals won by 3 wickets (with 2 balls remaining)"

# Fill in the matchNo value for the row-3 (now-shifted) data, which already
# carried the rest of its values over from the column insert
$ws.Cells.Item(3, 1).Value = "56th"
